$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("C4").Value = 100
$wsSummary.Range("A5").Value = 0
$wsSummary.Range("C5").Value = 0
$wsSummary.Range("D5").Select()

# ---------------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$wsRepay.Range("H3").Value = 42.74
$wsRepay.Range("K3").Value = 142.74
$wsRepay.Range("L3").Value = 142.74
$wsRepay.Range("N3").Value = 142.74

$wsRepay.Range("H4").Value = 3.25
$wsRepay.Range("K4").Value = 9903.25
$wsRepay.Range("L4").Value = 9903.25

# Column P was a duplicate "heading" placeholder column - remove it and
# shift the "Outstanding" column (was Q) one to the left, into P.
$wsRepay.Columns("P").Delete()

# Row 2's empty placeholder cell ends up one column further left (O2)
# rather than directly shifted into P2.
$wsRepay.Range("P2").Copy($wsRepay.Range("O2"))
$wsRepay.Range("P2").Clear()

$wsRepay.Range("H3:H4").Select()

# ---------------------------------------------------------------------------
# Sheet: Transactions
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value = 4586
$wsTxn.Range("E2").Value = 10045.99
$wsTxn.Range("G2").Value = 45.99

$wsTxn.Range("A3").Value = 4584
$wsTxn.Range("H3").Value = 100
$wsTxn.Range("I3").Value = 0

$wsTxn.Range("A4").Value = 4583

# Columns K:L are no longer used on this sheet - clear them out entirely.
$wsTxn.Range("K2:L4").Clear()

$wsTxn.Range("D4").Select()

# Activate the Transactions sheet last so it becomes the active tab
# (activeTab=4 / tabSelected on this sheet, removed from NewLoanInput).
$wsTxn.Activate()
